$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<name>_old" -> "<name>_FV2310", "<name>_new" ->
#    "<name>_FV2404". The "diff" header (K1) stays as-is.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# ---------------------------------------------------------------------------
# 2. Turn A1:U75 into a real table ("Table1") while preserving the existing
#    header-row cell formatting exactly (no new dxf / style should appear).
#    Creating a ListObject on a header row that already carries an explicit
#    style causes Excel to capture that style as the table's headerRowDxfId
#    override, which is NOT part of the target formatting. We sidestep this
#    by stashing the header formatting on a scratch row, clearing the header
#    row format before the table is created (so the header picks up no
#    explicit style / no dxf), then restoring the original formatting
#    afterwards via copy/paste (which reuses the existing style, unlike
#    reassigning Font/Interior properties one by one).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

$headerRange.Copy()
$scratchRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U75")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$scratchRange.ClearContents()
$scratchRange.ClearFormats()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (View > Freeze Panes > Freeze Top Row).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
